# Normalise the "sex" column (H) in Ark1 so every row uses the same
# lower-case coding: "M" -> "m", "K" -> "f", "F" -> "f".
# (Commit message: "endret på kjønn i exel slik at alle sier det samme")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("H2").Value  = "m"
$ws.Range("H3").Value  = "m"
$ws.Range("H4").Value  = "m"
$ws.Range("H5").Value  = "m"
$ws.Range("H6").Value  = "f"
$ws.Range("H7").Value  = "f"
$ws.Range("H8").Value  = "m"
$ws.Range("H9").Value  = "m"
$ws.Range("H10").Value = "m"
$ws.Range("H11").Value = "m"
$ws.Range("H12").Value = "f"

# Move the active selection to H12 (matches the saved cursor position).
$ws.Range("H12").Select() | Out-Null
